$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "Defect Size (mm)" (column D) to hold the new
# "Wheel Size" field. This shifts the existing D:H columns to E:I.
$ws.Columns("D:D").Insert()

# New column header
$ws.Range("D1").Value = "Wheel Size"

# Reorder the three "Axle 4 / Wheel 6/7/8" measurement rows so that the
# Axle 3 / Wheel 6 reading (defect 30.2994) now comes first (row 2),
# followed by Axle 4 / Wheel 8 (row 3) and Axle 4 / Wheel 7 (row 4).
# Rows 5 and 6 (Axle 3/Wheel 5 and Axle 2/Wheel 4) are unaffected.
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 6
$ws.Range("E2").Value = 30.2994

$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 8
$ws.Range("E3").Value = 30.1317

$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 7
$ws.Range("E4").Value = 30.064

# Populate the new "Wheel Size" column. Only the Axle 3 / Wheel 6 row
# carries a measured value; the rest report "0 mm".
$ws.Range("D2").Value = "73.25 in"
$ws.Range("D3").Value = "0 mm"
$ws.Range("D4").Value = "0 mm"
$ws.Range("D5").Value = "0 mm"
$ws.Range("D6").Value = "0 mm"
